# [FEATURE] Nuevos casos de Bloqueo y Desbloqueo (BYD 01-05)
# Adds a new user/branch pair (F00076 / 076) as the last row of the
# "Users" sheet, and leaves the sheet scrolled/selected on that new row,
# mirroring the state Excel saved after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Activate()

# New BYD case: user F00076 tied to branch 076
$ws.Range("A20").Value = "F00076"
$ws.Range("C20").Value = "076"

# Match the view state left behind in the saved workbook: scrolled so
# row 4 is at the top, with the newly-added C20 cell selected.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C20").Select()
